# Fruta / hortaliza, semanal
# Insert a new weekly record above row 5, shifting the existing rows 5-9 down to 6-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (pushes existing row 5 and below down by one)
$ws.Rows("5:5").Insert()

# Populate the new row 5 with the weekly data point.
# Columns A,B,C,E,F,G,H,I,J,K,L,N,Q,R,T mirror the record that used to sit in row 5
# (now shifted to row 6); D, M, O, P, S carry the new values for this week.
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Femacal de La Calera"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44589
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101004
$ws.Range("J5").Value = "Frambuesa"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 6000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("Q5").Value = "$/bandeja 2 kilos"
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 2
